# Bugfixed evaluation and simulated rt_data for components
#
# The "date" column (A2:A39) currently holds text labels like "1987Q4",
# "2024Q4", ... as shared strings. Replace them with real Excel date
# serial numbers (Dec 31 of each year) and format the column with a
# custom date/time number format, matching the fixed rt_data series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel date serial numbers (1899-12-30 epoch) for 1987-12-31 .. 2024-12-31 (Q4 of each year)
$dates = 32142,32508,32873,33238,33603,33969,34334,34699,35064,35430,35795,36160,36525,36891,37256,37621,37986,38352,38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Apply the custom date/time format to the whole date column range
$ws.Range("A2:A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
